$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 79.47369
$ws.Range("I9").Value = 70.07143000000001
$ws.Range("J9").Value = 105.8
$ws.Range("K9").Value = 70.07143000000001
$ws.Range("L9").Value = 105.8
$ws.Range("M9").Value = 98.92856999999999
$ws.Range("N9").Value = -443.8
$ws.Range("H32").Value = 27274118
$ws.Range("I32").Value = 50000976
$ws.Range("K32").Value = 50000976
$ws.Range("M32").Value = -50000650
$ws.Range("H40").Value = 4663.6665
$ws.Range("I40").Value = 4663.6665
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4663.6665
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4488.6665
$ws.Range("H76").Value = 5037
$ws.Range("J76").Value = 5599
$ws.Range("L76").Value = 5599
$ws.Range("N76").Value = -6229
$ws.Range("H79").Value = 5037
$ws.Range("J79").Value = 5599
$ws.Range("L79").Value = 5599
$ws.Range("N79").Value = -7783
$ws.Range("H107").Value = 1743.9286
$ws.Range("I107").Value = 1889.5
$ws.Range("J107").Value = 1549.8334
$ws.Range("K107").Value = 1889.5
$ws.Range("L107").Value = 1549.8334
$ws.Range("M107").Value = 30.5
$ws.Range("N107").Value = -5389.8334
$ws.Range("H132").Value = 5892.4067
$ws.Range("I132").Value = 3598.2856
$ws.Range("K132").Value = 10794.8568
$ws.Range("M132").Value = -8264.856800000001
$ws.Range("H136").Value = 150349.5
$ws.Range("J136").Value = 156264.4
$ws.Range("L136").Value = 156264.4
$ws.Range("N136").Value = -166464.4
$ws.Range("H137").Value = 5662.6
$ws.Range("I137").Value = 7501.2666
$ws.Range("K137").Value = 22503.7998
$ws.Range("M137").Value = -19953.7998
$ws.Range("N40").ClearContents()

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 29799.2
$ws.Range("J55").Value = 29799.2
$ws.Range("L55").Value = 29799.2
$ws.Range("N55").Value = -30429.2
$ws.Range("H62").Value = 40244
$ws.Range("J62").Value = 40244
$ws.Range("L62").Value = 40244
$ws.Range("N62").Value = -41492
$ws.Range("H65").Value = 40244
$ws.Range("J65").Value = 40244
$ws.Range("L65").Value = 120732
$ws.Range("N65").Value = -126972
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("H97").Value = 41686.64
$ws.Range("I97").Value = 2016.6666
$ws.Range("J97").Value = 143695.14
$ws.Range("K97").Value = 2016.6666
$ws.Range("L97").Value = 143695.14
$ws.Range("M97").Value = -1520.6666
$ws.Range("N97").Value = -144687.14
$ws.Range("N68").ClearContents()
$ws.Range("N71").ClearContents()

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 926.0323
$ws.Range("I94").Value = 530.1667
$ws.Range("J94").Value = 1474.1538
$ws.Range("K94").Value = 530.1667
$ws.Range("L94").Value = 1474.1538
$ws.Range("M94").Value = -79.16669999999999
$ws.Range("N94").Value = -2376.1538
$ws.Range("H134").Value = 2825
$ws.Range("I134").Value = 2322.7273
$ws.Range("K134").Value = 6968.1819
$ws.Range("M134").Value = -4433.1819

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 183
$ws.Range("I22").Value = 264
$ws.Range("J22").Value = 129
$ws.Range("K22").Value = 264
$ws.Range("L22").Value = 129
$ws.Range("M22").Value = 86
$ws.Range("N22").Value = -829
$ws.Range("H74").Value = 77999.664
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 77999.664
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 77999.664
$ws.Range("N74").Value = -79747.664
$ws.Range("H77").Value = 77999.664
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 77999.664
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 233998.992
$ws.Range("N77").Value = -242734.992
$ws.Range("H99").Value = 10964.393
$ws.Range("I99").Value = 7296.6875
$ws.Range("K99").Value = 7296.6875
$ws.Range("M99").Value = -5798.6875
$ws.Range("H122").Value = 13824.5
$ws.Range("I122").Value = 13824.5
$ws.Range("K122").Value = 41473.5
$ws.Range("M122").Value = -39023.5
$ws.Range("H126").Value = 10964.393
$ws.Range("I126").Value = 7296.6875
$ws.Range("K126").Value = 21890.0625
$ws.Range("M126").Value = -19420.0625
$ws.Range("H132").Value = 7463.0513
$ws.Range("I132").Value = 2245.84
$ws.Range("K132").Value = 6737.52
$ws.Range("M132").Value = -4207.52
$ws.Range("H134").Value = 5323.5713
$ws.Range("I134").Value = 5274.9375
$ws.Range("J134").Value = 5479.2
$ws.Range("K134").Value = 15824.8125
$ws.Range("L134").Value = 16437.6
$ws.Range("M134").Value = -13289.8125
$ws.Range("N134").Value = -21507.6
$ws.Range("M74").ClearContents()
$ws.Range("M77").ClearContents()

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 5147.6665
$ws.Range("J36").Value = 7499.5
$ws.Range("L36").Value = 22498.5
$ws.Range("N36").Value = -22836.5
$ws.Range("H64").Value = 1640.25
$ws.Range("I64").Value = 1640.25
$ws.Range("K64").Value = 4920.75
$ws.Range("M64").Value = -4650.75
$ws.Range("H67").Value = 1640.25
$ws.Range("I67").Value = 1640.25
$ws.Range("K67").Value = 4920.75
$ws.Range("M67").Value = -3984.75

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 43246.25
$ws.Range("J69").Value = 43246.25
$ws.Range("L69").Value = 43246.25
$ws.Range("N69").Value = -44744.25
$ws.Range("H72").Value = 43246.25
$ws.Range("J72").Value = 43246.25
$ws.Range("L72").Value = 129738.75
$ws.Range("N72").Value = -137226.75
$ws.Range("H113").Value = 1862.7858
$ws.Range("I113").Value = 1886.1428
$ws.Range("J113").Value = 1839.4286
$ws.Range("K113").Value = 1886.1428
$ws.Range("L113").Value = 1839.4286
$ws.Range("M113").Value = 283.8571999999999
$ws.Range("N113").Value = -6179.4286
$ws.Range("H122").Value = 5615.4
$ws.Range("I122").Value = 5529.1665
$ws.Range("J122").Value = 5744.75
$ws.Range("K122").Value = 16587.4995
$ws.Range("L122").Value = 17234.25
$ws.Range("M122").Value = -14137.4995
$ws.Range("N122").Value = -22134.25
$ws.Range("H132").Value = 3331.2058
$ws.Range("I132").Value = 2398.8076
$ws.Range("K132").Value = 7196.4228
$ws.Range("M132").Value = -4666.4228

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 10705.134
$ws.Range("I136").Value = 12189.75
$ws.Range("K136").Value = 36569.25
$ws.Range("M136").Value = -34019.25

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3505.842
$ws.Range("I132").Value = 3506.018
$ws.Range("J132").Value = 3501
$ws.Range("K132").Value = 10518.054
$ws.Range("L132").Value = 10503
$ws.Range("M132").Value = -7988.054
$ws.Range("N132").Value = -15563
$ws.Range("H135").Value = 76577.5
$ws.Range("J135").Value = 76577.5
$ws.Range("L135").Value = 76577.5
$ws.Range("N135").Value = -86717.5
$ws.Range("H136").Value = 4078.5715
$ws.Range("I136").Value = 4255.778
$ws.Range("K136").Value = 12767.334
$ws.Range("M136").Value = -10217.334
